$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted ahead of the existing row 303, shifting
# every subsequent record down by one row (old row 303 -> new row 304, ...,
# old row 407 -> new row 408).
$ws.Rows.Item(303).Insert()

# Populate the newly inserted row 303 with the new weekly entry. Excel's
# Insert() already copied formatting from the row above (e.g. the date
# style on column D), so we only need to set the values.
$ws.Cells.Item(303, 1).Value = 1
$ws.Cells.Item(303, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(303, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(303, 4).Value = 45119
$ws.Cells.Item(303, 5).Value = 15
$ws.Cells.Item(303, 6).Value = "Fruta"
$ws.Cells.Item(303, 7).Value = 100108
$ws.Cells.Item(303, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(303, 9).Value = 100108006
$ws.Cells.Item(303, 10).Value = "Plátano"
$ws.Cells.Item(303, 11).Value = "Sin especificar"
$ws.Cells.Item(303, 12).Value = "Pintón"
$ws.Cells.Item(303, 13).Value = 400
$ws.Cells.Item(303, 14).Value = 14000
$ws.Cells.Item(303, 15).Value = 15000
$ws.Cells.Item(303, 16).Value = 14500
$ws.Cells.Item(303, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(303, 18).Value = "Ecuador"
$ws.Cells.Item(303, 19).Value = 725
$ws.Cells.Item(303, 20).Value = 20
